# Add new columns I (I0) and J (IF) to the sheet, matching the style of
# the existing header row and plain numeric data cells.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Headers in row 1
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Copy the formatting (style) of the existing H1 header cell onto the two
# new header cells so they match the bold/bordered/centered header style.
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)

# Data rows 2-10 for columns I and J
$dataI = @(9, 8, 7, 6, 8, 7, 8, 8, 8)
$dataJ = @(9, 9, 7, 8, 9, 9, 8, 8, 8)

for ($i = 0; $i -lt $dataI.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 9).Value = $dataI[$i]
    $ws.Cells.Item($row, 10).Value = $dataJ[$i]
}
